$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Dilution") ahead of the existing data columns,
# shifting the former D:J (Temp2..Pyr_Rep3) right to E:K.
$ws.Columns("D").Insert()

# Header + constant dilution factor for each data row.
$ws.Range("D1").Value = "Dilution"
$ws.Range("D2:D7").Value = 200

# Update the sheet view: new zoom level and selection.
$excel.ActiveWindow.Zoom = 139
$ws.Range("F1:H8").Select() | Out-Null
